$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -6004.444595981147
$ws.Range("C2").Value = 13573.29778373467
$ws.Range("D2").Value = -6950.007831991451
$ws.Range("E2").Value = -618.8453557620749
$ws.Range("F2").Value = 56.36079380897991
$ws.Range("G2").Value = 45.21879701406675
$ws.Range("H2").Value = 56.68861583259826
$ws.Range("I2").Value = 45.36862495871259
$ws.Range("J2").Value = 56.42082461936974
$ws.Range("K2").Value = 45.8929585636667
$ws.Range("L2").Value = 53.57777799425878
$ws.Range("O2").Value = 45.60957459466476
$ws.Range("P2").Value = 54.92858261713064
$ws.Range("R2").Value = 7.698600306701701
$ws.Range("S2").Value = -17.12936487437229
$ws.Range("T2").Value = 9.430764567670584
$ws.Range("X2").Value = -120.9519315019875
$ws.Range("Y2").Value = -157.919038527756
$ws.Range("Z2").Value = -102.4455186539984
$ws.Range("AE2").Value = -7.698600306701701
$ws.Range("AF2").Value = 9.430764567670584
$ws.Range("AG2").Value = 7.698600306701701
$ws.Range("AH2").Value = -17.12936487437229
$ws.Range("AI2").Value = 9.430764567670584
$ws.Range("AJ2").Value = 7.698600306701701
$ws.Range("AK2").Value = -9.430764567670584
$ws.Range("AL2").Value = 18.48355351288429
$ws.Range("AM2").Value = -27.73675993687877
$ws.Range("AN2").Value = -120.9519315019875
$ws.Range("AO2").Value = -157.919038527756
$ws.Range("AP2").Value = -102.4455186539984
$ws.Range("AQ2").Value = -18.48355351288429
$ws.Range("AR2").Value = 27.73675993687877
$ws.Range("AS2").Value = 56.36079380897991
$ws.Range("AT2").Value = 56.36079380897991
$ws.Range("AU2").Value = 56.68861583259837
$ws.Range("AV2").Value = 56.68861583259837
$ws.Range("AW2").Value = 56.68861583259826
$ws.Range("AX2").Value = 56.42082461936974
$ws.Range("AY2").Value = 56.42082461936974
$ws.Range("AZ2").Value = 45.21879701406675
$ws.Range("BA2").Value = 45.21879701406675
$ws.Range("BB2").Value = 45.36862495871259
$ws.Range("BC2").Value = 44.97511407476708
$ws.Range("BD2").Value = 45.68985904213145
$ws.Range("BE2").Value = 45.8929585636667
$ws.Range("BF2").Value = 45.89295856366664
$ws.Range("BG2").Value = 53.57777799425878
$ws.Range("BJ2").Value = 45.60957459466476
$ws.Range("BK2").Value = 54.92858261713064
